$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: set "FECHA REAL DE CIERRE" (E4) and mark STATUS (F4) as "Cerrada"
$ws.Range("E4").Value = 42387
$ws.Range("F4").Value = "Cerrada"

# Update the active selection/cell to E5 (matches recorded view state)
$ws.Range("E5").Select()
